# Auto-generated script applying scraped market-price / profit updates
# to the Gungnir_Profits leve-crafting tables across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6081706
$ws.Range("I86").Value = 51100
$ws.Range("J86").Value = 10102110
$ws.Range("K86").Value = 51100
$ws.Range("L86").Value = 10102110
$ws.Range("M86").Value = -49977
$ws.Range("N86").Value = -10104356
$ws.Range("H89").Value = 6081706
$ws.Range("I89").Value = 51100
$ws.Range("J89").Value = 10102110
$ws.Range("K89").Value = 255500
$ws.Range("L89").Value = 50510550
$ws.Range("M89").Value = -249884
$ws.Range("N89").Value = -50521782
$ws.Range("H107").Value = 208.6875
$ws.Range("I107").Value = 175.75
$ws.Range("J107").Value = 307.5
$ws.Range("K107").Value = 175.75
$ws.Range("L107").Value = 307.5
$ws.Range("M107").Value = 1744.25
$ws.Range("N107").Value = -4147.5
$ws.Range("H132").Value = 8202236
$ws.Range("I132").Value = 8776875
$ws.Range("J132").Value = 13625
$ws.Range("K132").Value = 26330625
$ws.Range("L132").Value = 40875
$ws.Range("M132").Value = -26328095
$ws.Range("N132").Value = -45935
$ws.Range("H138").Value = 2778.82
$ws.Range("I138").Value = 1054.4166
$ws.Range("J138").Value = 3748.7969
$ws.Range("K138").Value = 3163.2498
$ws.Range("L138").Value = 11246.3907
$ws.Range("M138").Value = 1976.7502
$ws.Range("N138").Value = -21526.3907

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 45487220
$ws.Range("I32").Value = 76946920
$ws.Range("J32").Value = 45423.777
$ws.Range("K32").Value = 76946920
$ws.Range("L32").Value = 45423.777
$ws.Range("M32").Value = -76946633
$ws.Range("N32").Value = -45997.777
$ws.Range("H61").Value = 1459.3231
$ws.Range("I61").Value = 1366.84
$ws.Range("J61").Value = 1767.6
$ws.Range("K61").Value = 1366.84
$ws.Range("L61").Value = 1767.6
$ws.Range("M61").Value = -1154.84
$ws.Range("N61").Value = -2191.6
$ws.Range("H74").Value = 2490.8125
$ws.Range("I74").Value = 2717.7083
$ws.Range("J74").Value = 1810.125
$ws.Range("K74").Value = 2717.7083
$ws.Range("L74").Value = 1810.125
$ws.Range("M74").Value = -1843.7083
$ws.Range("N74").Value = -3558.125
$ws.Range("H77").Value = 2490.8125
$ws.Range("I77").Value = 2717.7083
$ws.Range("J77").Value = 1810.125
$ws.Range("K77").Value = 13588.5415
$ws.Range("L77").Value = 9050.625
$ws.Range("M77").Value = -9220.541499999999
$ws.Range("N77").Value = -17786.625
$ws.Range("H122").Value = 9317.25
$ws.Range("I122").Value = 16884.75
$ws.Range("J122").Value = 1749.75
$ws.Range("K122").Value = 50654.25
$ws.Range("L122").Value = 5249.25
$ws.Range("M122").Value = -48204.25
$ws.Range("N122").Value = -10149.25
$ws.Range("H132").Value = 20836524
$ws.Range("I132").Value = 27027890
$ws.Range("J132").Value = 11023.546
$ws.Range("K132").Value = 81083670
$ws.Range("L132").Value = 33070.638
$ws.Range("M132").Value = -81081140
$ws.Range("N132").Value = -38130.638
$ws.Range("H136").Value = 1459.3231
$ws.Range("I136").Value = 1366.84
$ws.Range("J136").Value = 1767.6
$ws.Range("K136").Value = 4100.52
$ws.Range("L136").Value = 5302.799999999999
$ws.Range("M136").Value = -1550.52
$ws.Range("N136").Value = -10402.8

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2511.6206
$ws.Range("I20").Value = 2383.762
$ws.Range("J20").Value = 2847.25
$ws.Range("K20").Value = 2383.762
$ws.Range("L20").Value = 2847.25
$ws.Range("M20").Value = -2136.762
$ws.Range("N20").Value = -3341.25
$ws.Range("H87").Value = 33750
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 33750
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 33750
$ws.Range("N87").Value = -36246
$ws.Range("H90").Value = 33750
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 33750
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 101250
$ws.Range("N90").Value = -113730
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = 0

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 16673433
$ws.Range("I132").Value = 1197.2
$ws.Range("J132").Value = 33345668
$ws.Range("K132").Value = 3591.6
$ws.Range("L132").Value = 100037004
$ws.Range("M132").Value = -1061.6
$ws.Range("N132").Value = -100042064
$ws.Range("H134").Value = 2723.9092
$ws.Range("I134").Value = 3021.8
$ws.Range("J134").Value = 2085.5715
$ws.Range("K134").Value = 9065.400000000001
$ws.Range("L134").Value = 6256.7145
$ws.Range("M134").Value = -6530.400000000001
$ws.Range("N134").Value = -11326.7145
$ws.Range("H138").Value = 52410
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 52410
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 52410
$ws.Range("N138").Value = -62690

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2191
$ws.Range("I58").Value = 1152.5
$ws.Range("J58").Value = 2300.3157
$ws.Range("K58").Value = 3457.5
$ws.Range("L58").Value = 6900.9471
$ws.Range("M58").Value = -3329.5
$ws.Range("N58").Value = -7156.9471
$ws.Range("H131").Value = 762.89
$ws.Range("I131").Value = 540
$ws.Range("J131").Value = 774.62103
$ws.Range("K131").Value = 1620
$ws.Range("L131").Value = 2323.86309
$ws.Range("M131").Value = 3420
$ws.Range("N131").Value = -12403.86309

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3000280.5
$ws.Range("I7").Value = 3333334
$ws.Range("J7").Value = 2500700
$ws.Range("K7").Value = 3333334
$ws.Range("L7").Value = 2500700
$ws.Range("M7").Value = -3333222
$ws.Range("N7").Value = -2500924
$ws.Range("H8").Value = 3000280.5
$ws.Range("I8").Value = 3333334
$ws.Range("J8").Value = 2500700
$ws.Range("K8").Value = 3333334
$ws.Range("L8").Value = 2500700
$ws.Range("M8").Value = -3333195
$ws.Range("N8").Value = -2500978
$ws.Range("H70").Value = 3767.1667
$ws.Range("I70").Value = 3724.4167
$ws.Range("J70").Value = 3852.6667
$ws.Range("K70").Value = 3724.4167
$ws.Range("L70").Value = 3852.6667
$ws.Range("M70").Value = -3454.4167
$ws.Range("N70").Value = -4392.6667
$ws.Range("H73").Value = 3767.1667
$ws.Range("I73").Value = 3724.4167
$ws.Range("J73").Value = 3852.6667
$ws.Range("K73").Value = 3724.4167
$ws.Range("L73").Value = 3852.6667
$ws.Range("M73").Value = -2788.4167
$ws.Range("N73").Value = -5724.6667
$ws.Range("H132").Value = 5019.3657
$ws.Range("I132").Value = 1526.3448
$ws.Range("J132").Value = 13460.833
$ws.Range("K132").Value = 4579.0344
$ws.Range("L132").Value = 40382.499
$ws.Range("M132").Value = -2049.0344
$ws.Range("N132").Value = -45442.499

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 680.4783
$ws.Range("I100").Value = 632.55
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 1265.1
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -724.0999999999999
$ws.Range("N100").Value = -3082
